# Update "想去人数" (F column) counts across all four sheets to reflect
# newly-generated output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 390
$ws.Range("F5").Value = 1180
$ws.Range("F8").Value = 1224
$ws.Range("F9").Value = 1677
$ws.Range("F10").Value = 6192
$ws.Range("F12").Value = 1802
$ws.Range("F13").Value = 472
$ws.Range("F19").Value = 6533
$ws.Range("F26").Value = 12
$ws.Range("F29").Value = 1543
$ws.Range("F31").Value = 301
$ws.Range("F35").Value = 79
$ws.Range("F36").Value = 3885

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 335
$ws.Range("F5").Value = 198
$ws.Range("F8").Value = 439

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 2252
$ws.Range("F4").Value = 655
$ws.Range("F5").Value = 237

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 2252
$ws.Range("F4").Value = 655
$ws.Range("F5").Value = 390
$ws.Range("F7").Value = 1180
$ws.Range("F11").Value = 335
$ws.Range("F12").Value = 1224
$ws.Range("F13").Value = 237
$ws.Range("F14").Value = 1677
$ws.Range("F15").Value = 6192
$ws.Range("F16").Value = 1802
$ws.Range("F19").Value = 472
$ws.Range("F24").Value = 6533
$ws.Range("F31").Value = 12
$ws.Range("F34").Value = 1544
$ws.Range("F37").Value = 301
$ws.Range("F44").Value = 79
$ws.Range("F46").Value = 3885
